$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

$ws.Range("G2").Value = 0.3776785511166963
$ws.Range("H2").Value = 0.968

$ws.Range("G3").Value = 0.3776785511166963
$ws.Range("H3").Value = 0.968

$ws.Range("G4").Value = 0.3776785511166963
$ws.Range("H4").Value = 0.968
